# Preparation for transport:
#  - Independent num/denum conversion
#  - Added some passenger convs
#  - CAP2ACT is now entity dependent
#
# Concretely: insert a new parameter row ("capacity_to_activity") right
# above the existing "co2_factor" row (row 10) in the conv_chp_oil block,
# shifting every row below it down by one, then fix up the sheet-level
# bookkeeping (used range, autofilter, filter-database defined name,
# active selection) that Excel keeps in sync when a row is inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at row 10 (pushes old rows 10.. down to 11..)
$ws.Rows.Item(10).Insert()

# 2. Populate the newly inserted row with the capacity_to_activity entry
$ws.Range("A10").Value = "CHE"
$ws.Range("B10").Value = "conv_chp_oil"
$ws.Range("C10").Value = "capacity_to_activity"
$ws.Range("D10").Value = "constant"
$ws.Range("G10").Value = 0.001
$ws.Range("H10").Value = "GW/TWh"

# 3. The data block now runs one row further, so the autofilter /
#    filter-database range needs to grow from L849 to L850 as well.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A5:L850").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$5:`$L`$850"
    }
}

# 4. Excel leaves the active selection on the cell that was pushed down
#    along with the insert (D9 -> D10).
$ws.Range("D10").Select()
